# Updated cryptos list values (Price and Volume(1h) columns) per source diff.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "30.614.08"
$ws.Range("E2").Value = "  +0.55%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.881.28"
$ws.Range("E3").Value = "  -0.53%  "
$ws.Range("E4").Value = "  +0.10%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "239.53"
$ws.Range("E5").Value = "  -0.27%  "
$ws.Range("E6").Value = "  +0.09%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4823"
$ws.Range("E7").Value = "  -0.49%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.2833"
$ws.Range("E8").Value = "  -2.16%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.06528"
$ws.Range("E9").Value = "  -1.57%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "1.962.36"
$ws.Range("E10").Value = "  +3.67%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07529"
$ws.Range("E11").Value = "  +1.58%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "16.62"
$ws.Range("E12").Value = "  -2.16%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "5.117"
$ws.Range("E13").Value = "  -1.72%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "88.71"
$ws.Range("E14").Value = "  -0.65%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.6657"
$ws.Range("E15").Value = "  +0.07%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "30.560.81"
$ws.Range("E16").Value = "  +0.47%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "13.37"
$ws.Range("E17").Value = "  -1.40%  "
$ws.Range("E18").Value = "  +0.09%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.000007610"
$ws.Range("E19").Value = "  -2.33%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "230.93"
$ws.Range("E20").Value = "  +6.03%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "2.136.07"
$ws.Range("E21").Value = "  +0.09%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "5.294"
$ws.Range("E22").Value = "  -2.65%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "1.001"
$ws.Range("E23").Value = "  +0.13%  "
$ws.Range("E24").Value = "  -0.49%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "9.381"
$ws.Range("E25").Value = "  -0.77%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "167.75"
$ws.Range("E26").Value = "  +1.56%  "
$ws.Range("E27").Value = "  +0.20%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "1.946"
$ws.Range("E28").Value = "  -0.14%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.423"
$ws.Range("E29").Value = "  -1.22%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.09584"
$ws.Range("E30").Value = "  +4.20%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "4.371"
$ws.Range("E31").Value = "  +1.06%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "4.042"
$ws.Range("E32").Value = "  -0.93%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.05035"
$ws.Range("E33").Value = "  -0.98%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.208"
$ws.Range("E34").Value = "  +4.23%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.7475"
$ws.Range("E35").Value = "  -0.40%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "2.700"
$ws.Range("E36").Value = "  -0.14%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.01852"
$ws.Range("E37").Value = "  -2.75%  "
$ws.Range("E38").Value = "  -0.44%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "2.095"
$ws.Range("E39").Value = "  -0.20%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.9174"
$ws.Range("E40").Value = "  -0.44%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "106.21"
$ws.Range("E41").Value = "  -1.25%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.4284"
$ws.Range("E42").Value = "  -1.79%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "5.794"
$ws.Range("E43").Value = "  -4.89%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "1.000"
$ws.Range("E44").Value = "  -0.12%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "7.416"
$ws.Range("E45").Value = "  -3.21%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.1286"
$ws.Range("E46").Value = "  -4.72%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "64.20"
$ws.Range("E47").Value = "  -3.06%  "
$ws.Range("E48").Value = "  -6.25%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "8.952"
$ws.Range("E49").Value = "  -0.09%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "33.89"
$ws.Range("E50").Value = "  -1.66%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.05647"
$ws.Range("E51").Value = "  -0.93%  "
